$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row before row 546; this shifts rows 546:599 down to 547:600
$ws.Rows.Item(546).Insert()

# Populate the newly inserted row 546 with the new data record
$ws.Cells.Item(546, 1).Value = 5
$ws.Cells.Item(546, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(546, 3).Value = "Maule"
$ws.Cells.Item(546, 4).Value = 45132
$ws.Cells.Item(546, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(546, 5).Value = 7
$ws.Cells.Item(546, 6).Value = 100112032
$ws.Cells.Item(546, 7).Value = "Zapallo italiano"
$ws.Cells.Item(546, 8).Value = "Sin especificar"
$ws.Cells.Item(546, 9).Value = "Primera"
$ws.Cells.Item(546, 10).Value = 300
$ws.Cells.Item(546, 11).Value = 14000
$ws.Cells.Item(546, 12).Value = 14000
$ws.Cells.Item(546, 13).Value = 14000
$ws.Cells.Item(546, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(546, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(546, 16).Value = 280
$ws.Cells.Item(546, 17).Value = 50
$ws.Cells.Item(546, 18).Value = "Hortaliza"
